$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1682847896440129
$ws.Range("C2").Value = 0.598705501618123
$ws.Range("J2").Value = 0.0151024811218986
$ws.Range("P2").Value = 0.1143473570658037
$ws.Range("S2").Value = 0.1035598705501618
$ws.Range("B3").Value = 0.01208981001727116
$ws.Range("C3").Value = 0.03281519861830743
$ws.Range("J3").Value = 0.0310880829015544
$ws.Range("P3").Value = 0.7633851468048359
$ws.Range("S3").Value = 0.1606217616580311
$ws.Range("J4").Value = 0.0736196319018405
$ws.Range("O4").Value = 0.006134969325153374
$ws.Range("P4").Value = 0.6748466257668712
$ws.Range("S4").Value = 0.245398773006135
$ws.Range("B6").Value = 0.06732117812061711
$ws.Range("D6").Value = 0.01262272089761571
$ws.Range("E6").Value = 0.002805049088359046
$ws.Range("F6").Value = 0.06872370266479663
$ws.Range("J6").Value = 0.2356241234221599
$ws.Range("O6").Value = 0.01122019635343618
$ws.Range("Q6").Value = 0.1598877980364656
$ws.Range("R6").Value = 0.06451612903225806
$ws.Range("S6").Value = 0.3772791023842917
$ws.Range("B7").Value = 0.1234782608695652
$ws.Range("D7").Value = 0.02260869565217391
$ws.Range("E7").Value = 0.001739130434782609
$ws.Range("F7").Value = 0.05391304347826087
$ws.Range("J7").Value = 0.1252173913043478
$ws.Range("O7").Value = 0.03304347826086956
$ws.Range("Q7").Value = 0.1791304347826087
$ws.Range("R7").Value = 0.06956521739130435
$ws.Range("S7").Value = 0.391304347826087
$ws.Range("B8").Value = 0.105379513633014
$ws.Range("D8").Value = 0.02431834929992631
$ws.Range("F8").Value = 0.07737656595431099
$ws.Range("J8").Value = 0.1127487103905674
$ws.Range("O8").Value = 0.02210759027266028
$ws.Range("Q8").Value = 0.182756079587325
$ws.Range("R8").Value = 0.1024318349299926
$ws.Range("S8").Value = 0.3728813559322034
$ws.Range("B9").Value = 0.1096774193548387
$ws.Range("D9").Value = 0.02419354838709677
$ws.Range("E9").Value = 0.001612903225806452
$ws.Range("F9").Value = 0.07258064516129033
$ws.Range("J9").Value = 0.09516129032258064
$ws.Range("O9").Value = 0.02419354838709677
$ws.Range("Q9").Value = 0.1532258064516129
$ws.Range("R9").Value = 0.09677419354838709
$ws.Range("S9").Value = 0.4225806451612903
$ws.Range("B10").Value = 0.1055681531872008
$ws.Range("D10").Value = 0.02443940539178634
$ws.Range("E10").Value = 0.002015621063240111
$ws.Range("F10").Value = 0.07180650037792895
$ws.Range("J10").Value = 0.1118669690098262
$ws.Range("O10").Value = 0.02091206853111615
$ws.Range("Q10").Value = 0.2091206853111615
$ws.Range("R10").Value = 0.08742756361803981
$ws.Range("S10").Value = 0.3668430335097002
$ws.Range("G11").Value = 0.1407089151450054
$ws.Range("J11").Value = 0.1160042964554243
$ws.Range("K11").Value = 0.2083780880773362
$ws.Range("L11").Value = 0.5263157894736842
$ws.Range("S11").Value = 0.008592910848549946
$ws.Range("G12").Value = 0.708984375
$ws.Range("J12").Value = 0.208984375
$ws.Range("K12").Value = 0.00390625
$ws.Range("L12").Value = 0.02734375
$ws.Range("S12").Value = 0.05078125
$ws.Range("G13").Value = 0.6901408450704225
$ws.Range("J13").Value = 0.2605633802816901
$ws.Range("S13").Value = 0.04929577464788732
$ws.Range("F15").Value = 0.02624671916010499
$ws.Range("H15").Value = 0.1548556430446194
$ws.Range("I15").Value = 0.07086614173228346
$ws.Range("J15").Value = 0.3412073490813649
$ws.Range("K15").Value = 0.06692913385826772
$ws.Range("M15").Value = 0.01049868766404199
$ws.Range("O15").Value = 0.06430446194225722
$ws.Range("S15").Value = 0.2650918635170604
$ws.Range("F16").Value = 0.01399688958009331
$ws.Range("H16").Value = 0.1539657853810264
$ws.Range("I16").Value = 0.08087091757387248
$ws.Range("J16").Value = 0.4199066874027994
$ws.Range("K16").Value = 0.1010886469673406
$ws.Range("M16").Value = 0.02488335925349922
$ws.Range("O16").Value = 0.05909797822706065
$ws.Range("S16").Value = 0.1461897356143079
$ws.Range("F17").Value = 0.0180766449746927
$ws.Range("H17").Value = 0.1778741865509761
$ws.Range("I17").Value = 0.08821402747650037
$ws.Range("J17").Value = 0.4193781634128706
$ws.Range("K17").Value = 0.09038322487346348
$ws.Range("M17").Value = 0.02241503976861894
$ws.Range("N17").Value = 0.001446131597975416
$ws.Range("O17").Value = 0.06146059291395517
$ws.Range("S17").Value = 0.1207519884309472
$ws.Range("F18").Value = 0.02866242038216561
$ws.Range("H18").Value = 0.1942675159235669
$ws.Range("I18").Value = 0.07802547770700637
$ws.Range("J18").Value = 0.410828025477707
$ws.Range("K18").Value = 0.09076433121019108
$ws.Range("M18").Value = 0.0143312101910828
$ws.Range("N18").Value = 0.001592356687898089
$ws.Range("O18").Value = 0.07006369426751592
$ws.Range("S18").Value = 0.1114649681528662
$ws.Range("F19").Value = 0.01386036960985626
$ws.Range("H19").Value = 0.1996919917864476
$ws.Range("I19").Value = 0.08752566735112936
$ws.Range("J19").Value = 0.3788501026694045
$ws.Range("K19").Value = 0.1095995893223819
$ws.Range("M19").Value = 0.02130390143737166
$ws.Range("N19").Value = 0.0007700205338809035
$ws.Range("O19").Value = 0.07751540041067762
$ws.Range("S19").Value = 0.1108829568788501
